$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 531
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 184

# Remove row 4 entirely (was A4=1, B4=184)
$ws.Rows.Item(4).Delete()
